# Add a new "test_xlr_n_percent" column to the existing table (table_test_1),
# fill it with the character representation "n (p%)" for each of the 32 rows,
# and give the new column the same right/bottom aligned, General-number-format
# style used by the other "no explicit number format" text columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Grow the table by one column (this also extends ref / autoFilter / dimension).
$newCol = $lo.ListColumns.Add()
$ws.Range("K2").Value = "test_xlr_n_percent"

$values = @("1 (3%)","2 (6%)","3 (9%)","4 (12%)","5 (16%)","6 (19%)","7 (22%)","8 (25%)","9 (28%)","10 (31%)","11 (34%)","12 (38%)","13 (41%)","14 (44%)","15 (47%)","16 (50%)","17 (53%)","18 (56%)","19 (59%)","20 (62%)","21 (66%)","22 (69%)","23 (72%)","24 (75%)","25 (78%)","26 (81%)","27 (84%)","28 (88%)","29 (91%)","30 (94%)","31 (97%)","32 (100%)")

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 11).Value = $values[$i]
}

# Base the new column's formatting on an existing plain-text column (font +
# general number format), then flip the alignment to right/bottom to match
# the numeric-looking "n (p%)" strings.
$src = $ws.Range("G3")
$src.Copy()
$first = $ws.Range("K3")
$first.PasteSpecial(-4122)
$first.VerticalAlignment = -4107
$first.HorizontalAlignment = -4152

# Propagate the exact resulting style (a single new cellXf) to the rest of
# the column, including the trailing blank row used by the table.
$first.Copy()
$rest = $ws.Range("K4:K35")
$rest.PasteSpecial(-4122)

$excel.CutCopyMode = 0
